$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: Category / Subcategory swap places (B1 <-> C1 text) ---
$ws.Range("B1").Value = "Subcategory"
$ws.Range("C1").Value = "Category"

# Columns J (Size) and R (ItemMRP) hold numeric-looking text ("32", "2195.00")
# that Excel would otherwise auto-convert to a number. Pre-format as Text so
# the values stick as strings, then clear the formatting override back off
# (the source file doesn't carry any number-format styling on these cells).
$textRangeJ = $ws.Range("J2:J5")
$textRangeR = $ws.Range("R2:R5")
$textRangeJ.NumberFormat = "@"
$textRangeR.NumberFormat = "@"

# --- Row 2: update existing data row with new prefixed/suffixed values ---
$ws.Range("A2").Value = "TROUSERS"
$ws.Range("B2").Value = "SC1"
$ws.Range("C2").Value = "C1"
$ws.Range("D2").Value = "pb141141ds1cs2ss1"
$ws.Range("F2").Value = "pcs"
$ws.Range("G2").Value = 40
$ws.Range("I2").Value = "L.GREY"
$ws.Range("J2").Value = "32"
$ws.Range("L2").Value = "B1"
$ws.Range("N2").Value = "sup12"
$ws.Range("R2").Value = "2195.00"
$ws.Range("T2").Value = 0

# --- Row 3: new row (same item master data, different description/size count) ---
$ws.Range("A3").Value = "TROUSERS"
$ws.Range("B3").Value = "SC1"
$ws.Range("C3").Value = "C1"
$ws.Range("D3").Value = "pb141141ds1cs2ss1"
$ws.Range("F3").Value = "pcs"
$ws.Range("G3").Value = 39
$ws.Range("I3").Value = "L.GREY"
$ws.Range("J3").Value = "32"
$ws.Range("L3").Value = "B1"
$ws.Range("N3").Value = "sup12"
$ws.Range("R3").Value = "2195.00"
$ws.Range("T3").Value = 0

# --- Row 4: new row ---
$ws.Range("A4").Value = "TROUSERS"
$ws.Range("B4").Value = "SC1"
$ws.Range("C4").Value = "C1"
$ws.Range("D4").Value = "pb141141ds1cs2ss1"
$ws.Range("F4").Value = "pcs"
$ws.Range("G4").Value = 38
$ws.Range("I4").Value = "L.GREY"
$ws.Range("J4").Value = "32"
$ws.Range("L4").Value = "B1"
$ws.Range("N4").Value = "sup12"
$ws.Range("R4").Value = "2195.00"
$ws.Range("T4").Value = 0

# --- Row 5: new row ---
$ws.Range("A5").Value = "TROUSERS"
$ws.Range("B5").Value = "SC1"
$ws.Range("C5").Value = "C1"
$ws.Range("D5").Value = "pb141141ds1cs2ss1"
$ws.Range("F5").Value = "pcs"
$ws.Range("G5").Value = 36
$ws.Range("I5").Value = "L.GREY"
$ws.Range("J5").Value = "32"
$ws.Range("L5").Value = "B1"
$ws.Range("N5").Value = "sup12"
$ws.Range("R5").Value = "2195.00"
$ws.Range("T5").Value = 0

# Drop the temporary Text number format again now the values are locked in
# as strings, so the cells end up styled the same as the rest of the sheet.
$textRangeJ.ClearFormats()
$textRangeR.ClearFormats()
